# Applies the RECO_holdings.xlsx update:
#   1. Updates the "as of" date in the confidential disclosure note
#      (A41) from 2021-04-29 to 2021-04-30.
#   2. Refreshes the Weight (col D) and Percent Change (col E) values
#      for every holding row (2-38) with the latest model figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (locked cells); unprotect so the refreshed
# figures can be written, then re-protect to restore the original state.
$ws.Unprotect()

try {
    $newNote = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-30 for illustrative purposes only and are subject to change."
    $ws.Range("A41").Value = $newNote

    $ws.Range("D2").Value = 0.02987848333873181
    $ws.Range("E2").Value = 0.005500946884299696
    $ws.Range("D3").Value = 0.02979784303547808
    $ws.Range("E3").Value = -0.01501676064252355
    $ws.Range("D4").Value = 0.03068353915856559
    $ws.Range("E4").Value = -0.02729114527469889
    $ws.Range("D5").Value = 0.06680847042666034
    $ws.Range("E5").Value = -0.001120614407817189
    $ws.Range("D6").Value = 0.0153154989315641
    $ws.Range("E6").Value = 0.005202442886224734
    $ws.Range("D7").Value = 0.01647602644617029
    $ws.Range("E7").Value = -0.02719372021306421
    $ws.Range("D8").Value = 0.02710322516947455
    $ws.Range("E8").Value = 0.02100464402880142
    $ws.Range("D9").Value = 0.03445015341865809
    $ws.Range("E9").Value = -0.01128491620111738
    $ws.Range("D10").Value = 0.02938848283256475
    $ws.Range("E10").Value = 0.008781925343811503
    $ws.Range("D11").Value = 0.03154902709110026
    $ws.Range("E11").Value = -0.01379280895038015
    $ws.Range("D12").Value = 0.01301792389518454
    $ws.Range("E12").Value = -0.01426670609107039
    $ws.Range("D13").Value = 0.01421309402216704
    $ws.Range("E13").Value = -0.04322274881516597
    $ws.Range("D14").Value = 0.01577316633284661
    $ws.Range("E14").Value = -0.008711991800478103
    $ws.Range("D15").Value = 0.008804612394395365
    $ws.Range("E15").Value = -0.003235114103348868
    $ws.Range("D16").Value = 0.007697588470014798
    $ws.Range("E16").Value = -0.01575157515751568
    $ws.Range("D17").Value = 0.03157327691976367
    $ws.Range("E17").Value = -0.01723843659327529
    $ws.Range("D18").Value = 0.02710245533364396
    $ws.Range("E18").Value = -0.032573035463209
    $ws.Range("D19").Value = 0.03175572801161221
    $ws.Range("E19").Value = -0.01181818181818184
    $ws.Range("D20").Value = 0.0317085755669889
    $ws.Range("E20").Value = -0.01344420503171384
    $ws.Range("D21").Value = 0.04676540965947659
    $ws.Range("E21").Value = -0.008136170773162532
    $ws.Range("D22").Value = 0.03396611914017807
    $ws.Range("E22").Value = -0.01280562087429538
    $ws.Range("D23").Value = 0.03130729864029671
    $ws.Range("E23").Value = -0.005133091535009293
    $ws.Range("D24").Value = 0.02979611090485926
    $ws.Range("E24").Value = -0.007615393558888517
    $ws.Range("D25").Value = 0.01470578895374114
    $ws.Range("E25").Value = 0.02414605418138982
    $ws.Range("D26").Value = 0.01489978758304845
    $ws.Range("E26").Value = -0.004882585445245224
    $ws.Range("D27").Value = 0.0305363080559663
    $ws.Range("E27").Value = -0.0163616195230174
    $ws.Range("D28").Value = 0.02929533269706402
    $ws.Range("E28").Value = -0.01424291795869037
    $ws.Range("D29").Value = 0.02915868683713526
    $ws.Range("E29").Value = -0.001306878935487754
    $ws.Range("D30").Value = 0.02821178876551628
    $ws.Range("E30").Value = -0.00480264145279885
    $ws.Range("D31").Value = 0.03539262493425121
    $ws.Range("E31").Value = -0.02057129806359004
    $ws.Range("D32").Value = 0.03047279659994307
    $ws.Range("E32").Value = 0.001528414617201657
    $ws.Range("D33").Value = 0.02979842041235102
    $ws.Range("E33").Value = -0.01981528127623855
    $ws.Range("D34").Value = 0.03091275777712214
    $ws.Range("E34").Value = -0.006848462208940309
    $ws.Range("D35").Value = 0.02988618169703765
    $ws.Range("E35").Value = 0.000115915150110224
    $ws.Range("D36").Value = 0.02986770563710363
    $ws.Range("E36").Value = -0.008892325536439372
    $ws.Range("D37").Value = 0.03192971090932431
    $ws.Range("E37").Value = -0.0005786478927573002
    $ws.Range("E38").Value = -0.008651607523066485
}
finally {
    $ws.Protect()
}

Write-Output "RECO_holdings.xlsx: disclosure date and holdings figures updated."
